# Update SWOT Matrix and Zone to Win risk analysis
$p = $ppt.ActivePresentation

# --- Slide 2 ("Methodology"): rewrite the intro paragraph with the
#     proxy-data disclaimer and tighter wording. ---
$s2 = $p.Slides.Item(2)
$methodologyShape = $s2.Shapes.Item(2)
$methodologyRange = $methodologyShape.TextFrame.TextRange
$introPara = $methodologyRange.Paragraphs(1)
# Clear first, then set, so the host doesn't diff/split the run against the
# old text (it otherwise keeps the common prefix/suffix as separate runs).
$introPara.Text = ""
$introPara.Text = "Gray Associates Portfolio Analysis evaluates academic programs by plotting Market Score (student demand 40% + employment 40% + competition 20%) against Program Economics (SCH efficiency + cost structure). Programs are classified as Grow, Sustain, Transform, Evaluate, or Sunset Review. Important: FLC does not have a Gray Associates subscription; scores are proxy estimates based on FLC institutional data, not official Gray output."

# --- Slide 8 ("Key Findings"): update the four category call-outs and
#     append a new disclaimer paragraph. ---
$s8 = $p.Slides.Item(8)
$findingsShape = $s8.Shapes.Item(2)
$findingsRange = $findingsShape.TextFrame.TextRange

$growPara = $findingsRange.Paragraphs(1)
$growPara.Text = ""
$growPara.Text = "GROW programs (high market + strong economics): Business Admin, Psychology, Engineering, Health Sciences, CIS, Exercise Physiology, Accounting show strongest investment case."

$sustainPara = $findingsRange.Paragraphs(2)
$sustainPara.Text = ""
$sustainPara.Text = "SUSTAIN programs (solid market, needs efficiency): Environmental programs, Criminology, Biology, Sociology, Teacher Education maintain enrollment but need optimization."

$transformPara = $findingsRange.Paragraphs(3)
$transformPara.Text = ""
$transformPara.Text = "TRANSFORM programs (weak market, strong economics): English and Math generate significant SCH as foundational/service courses — low Market Score reflects major enrollment, not institutional value."

$evaluatePara = $findingsRange.Paragraphs(4)
$evaluatePara.Text = ""
$evaluatePara.Text = "EVALUATE/SUNSET programs: Political Science, Art & Design require strategic review. Note: NAIS is mission-critical and must not be evaluated on enrollment metrics alone."

# Append the new closing disclaimer paragraph (keeps the same `<a:pPr/>`
# paragraph-formatting pattern the other non-first paragraphs use).
$null = $findingsRange.InsertAfter("`rData source disclaimer: FLC does not have a Gray Associates subscription. Scores are proxy estimates based on FLC data, not official Gray Associates output.")

Write-Host "Updated methodology and key findings slides"
